# #20 Time Tracking.xlsx updated to 2024-02-26.
# Appends new time-tracking session rows (2024-02-19, 2024-02-20, 2024-02-25,
# 2024-02-26) below the existing data, plus extra blank padding rows, and
# moves the active-cell selection down to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Session {
    param([int]$Row, [string]$Date, [string]$StartTime, [string]$EndTime, [string]$Effort, [string]$Hashtag, [string]$Descriptor)

    $ws.Range("A$Row").Value = $Date
    $ws.Range("B$Row").Value = $StartTime
    $ws.Range("C$Row").Value = $EndTime
    $ws.Range("D$Row").Value = $Effort
    $ws.Range("E$Row").Value = $Hashtag
    if ($Descriptor -ne "") {
        $ws.Range("F$Row").Value = $Descriptor
    }
    # Leading apostrophe forces text storage ("False" as a label, not a bool).
    $ws.Range("G$Row").Value = "'False"
    $ws.Range("H$Row").Value = "'False"
    $ws.Range("I$Row").Formula = "=YEAR(A$Row)"
    $ws.Range("J$Row").Formula = "=MONTH(A$Row)"
}

# 2024-02-19 -- #studying / Books.
Set-Session 994 "2024-02-19" "11:15" "13:00" "1h 45m" "#studying" "Books."
Set-Session 995 "2024-02-19" "15:30" "18:00" "2h 30m" "#studying" "Books."
Set-Session 996 "2024-02-19" "20:15" "21:15" "1h 00m" "#studying" "Books."

# 2024-02-20 -- #studying / Books.
Set-Session 997 "2024-02-20" "08:45" "12:15" "3h 30m" "#studying" "Books."
Set-Session 998 "2024-02-20" "13:30" "14:00" "0h 30m" "#studying" "Books."
Set-Session 999 "2024-02-20" "15:30" "16:30" "1h 00m" "#studying" "Books."

# 2024-02-25 -- #studying / Books.
Set-Session 1000 "2024-02-25" "10:15" "13:00" "2h 45m" "#studying" "Books."
Set-Session 1001 "2024-02-25" "14:00" "19:45" "5h 45m" "#studying" "Books."

# 2024-02-26 -- #studying / Books., then #maintenance (no descriptor)
Set-Session 1002 "2024-02-26" "08:15" "12:45" "4h 30m" "#studying" "Books."
Set-Session 1003 "2024-02-26" "16:00" "17:00" "1h 00m" "#studying" "Books."
Set-Session 1004 "2024-02-26" "17:00" "18:00" "1h 00m" "#maintenance" ""
Set-Session 1005 "2024-02-26" "21:30" "22:30" "1h 00m" "#maintenance" ""

# Extend the blank trailing rows (previously ending at 1007) down to 1025,
# copying the formatting of an existing blank row so no new styles appear.
$ws.Range("A1006:J1006").Copy()
$ws.Range("A1008:J1025").PasteSpecial(-4122)

# Move the view: select the new bottom-most cell that matches the commit.
$ws.Range("D1007").Select()
